$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1) Update the "Date" metadata value on the Metadata sheet (B8)
# ---------------------------------------------------------------
$metaWs = $wb.Worksheets.Item("Metadata")
$metaWs.Range("B8").Value = "2025-10-30T16:59:08+00:00"

# ---------------------------------------------------------------
# 2) Add a new row (11) to the Elements sheet describing the
#    "Professionnel.PersonnePriseCharge" element.
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("Elements")

$newId   = "Professionnel.PersonnePriseCharge"
$newType = "https://interop.esante.gouv.fr/ig/mos/StructureDefinition/PersonnePriseCharge`n"
$newDesc = "Lien vers la classe PersonnePriseCharge"

# Clone the formatting (and "blank"/placeholder content) of the last
# data row (row 10) into the new row 11 so the new row keeps the
# exact same cell styles as the rest of the table.
$ws.Range("A10:AJ10").Copy()
$ws.Range("A11:AJ11").PasteSpecial(-4122)

# A couple of cells in row 10 hold plain numeric-looking text ("0")
# that Excel's automatic type inference would otherwise turn into a
# real number. Force them to stay text by using a leading quote,
# then restore the row's normal formatting on top (this keeps the
# text value while reapplying the original, shared cell style).
$ws.Range("F11").Value = "'0"
$ws.Range("AG11").Value = "'0"
$ws.Range("A10:AJ10").Copy()
$ws.Range("A11:AJ11").PasteSpecial(-4122)

# Fill in the values that differ from row 10.
$ws.Range("A11").Value = $newId
$ws.Range("B11").Value = $newId
$ws.Range("K11").Value = $newType
$ws.Range("L11").Value = $newDesc
$ws.Range("M11").Value = $newDesc
$ws.Range("AF11").Value = $newId

# The Type(s) cell contains an embedded newline which makes the
# engine auto-expand the row height; restore the row to its default
# (automatic) height so it matches the rest of the table.
$ws.Rows.Item(11).AutoFit()

# ---------------------------------------------------------------
# 3) Column K ("Type(s)") grew a bit wider to fit the new content.
# ---------------------------------------------------------------
$ws.Columns.Item(11).ColumnWidth = 62.66796875
